$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H33").Value = 19289412
$ws.Range("I33").Value = 107735.57
$ws.Range("K33").Value = 107735.57
$ws.Range("M33").Value = -107506.57

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H41").Value = 77279.46000000001
$ws.Range("J41").Value = 167030.67
$ws.Range("L41").Value = 167030.67
$ws.Range("N41").Value = -167910.67

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H46").Value = 11949.9
$ws.Range("J46").Value = 11949.9
$ws.Range("L46").Value = 35849.7
$ws.Range("N46").Value = -36087.7

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H60").Value = 11949.9
$ws.Range("J60").Value = 11949.9
$ws.Range("L60").Value = 35849.7
$ws.Range("N60").Value = -36817.7

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H62").Value = 1549715.2
$ws.Range("I62").Value = 2067201.4
$ws.Range("J62").Value = 256000
$ws.Range("K62").Value = 2067201.4
$ws.Range("L62").Value = 256000
$ws.Range("M62").Value = -2066577.4
$ws.Range("N62").Value = -257248

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H65").Value = 1549715.2
$ws.Range("I65").Value = 2067201.4
$ws.Range("J65").Value = 256000
$ws.Range("K65").Value = 10336007
$ws.Range("L65").Value = 1280000
$ws.Range("M65").Value = -10332887
$ws.Range("N65").Value = -1286240

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H74").Value = 4500
$ws.Range("J74").Value = 4500
$ws.Range("L74").Value = 4500
$ws.Range("N74").Value = -6372

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H77").Value = 4500
$ws.Range("J77").Value = 4500
$ws.Range("L77").Value = 22500
$ws.Range("N77").Value = -31860

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H98").Value = 965.12
$ws.Range("I98").Value = 920.8095
$ws.Range("K98").Value = 920.8095
$ws.Range("M98").Value = 577.1905

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H106").Value = 143007330
$ws.Range("I106").Value = 200009860
$ws.Range("J106").Value = 500999.5
$ws.Range("K106").Value = 200009860
$ws.Range("L106").Value = 500999.5
$ws.Range("M106").Value = -200009229
$ws.Range("N106").Value = -502261.5

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H122").Value = 965.12
$ws.Range("I122").Value = 920.8095
$ws.Range("K122").Value = 2762.4285
$ws.Range("M122").Value = -312.4285

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H129").Value = 1292.5555
$ws.Range("I129").Value = 1151.2
$ws.Range("J129").Value = 1999.3334
$ws.Range("K129").Value = 3453.6
$ws.Range("L129").Value = 5998.0002
$ws.Range("M129").Value = 1546.4
$ws.Range("N129").Value = -15998.0002

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H131").Value = 7260.615
$ws.Range("I131").Value = 2448.5
$ws.Range("K131").Value = 7345.5
$ws.Range("M131").Value = -2305.5

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H132").Value = 3276.2964
$ws.Range("I132").Value = 3201.761
$ws.Range("J132").Value = 3704.875
$ws.Range("K132").Value = 9605.282999999999
$ws.Range("L132").Value = 11114.625
$ws.Range("M132").Value = -7075.282999999999
$ws.Range("N132").Value = -16174.625

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H135").Value = 55556812

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H141").Value = 949.1724
$ws.Range("I141").Value = 949.1724
$ws.Range("K141").Value = 2847.5172
$ws.Range("M141").Value = 2332.4828

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H27").Value = 28750
$ws.Range("I27").Value = 50000
$ws.Range("K27").Value = 50000
$ws.Range("M27").Value = -49816

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H132").Value = 1791.421
$ws.Range("I132").Value = 1817.1333
$ws.Range("K132").Value = 5451.3999
$ws.Range("M132").Value = -2921.3999

$ws = $wb.Sheets.Item("BSM")
$ws.Range("H64").Value = 100000
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()

$ws = $wb.Sheets.Item("BSM")
$ws.Range("H67").Value = 100000
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()

$ws = $wb.Sheets.Item("BSM")
$ws.Range("H94").Value = 4246.625
$ws.Range("I94").Value = 3995.6667
$ws.Range("K94").Value = 3995.6667
$ws.Range("M94").Value = -3544.6667

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H7").Value = 146.23077
$ws.Range("I7").Value = 126.2
$ws.Range("J7").Value = 213
$ws.Range("K7").Value = 126.2
$ws.Range("L7").Value = 213
$ws.Range("M7").Value = -13.2
$ws.Range("N7").Value = -439

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H31").Value = 3003.5078
$ws.Range("I31").Value = 2301.6086
$ws.Range("K31").Value = 2301.6086
$ws.Range("M31").Value = -2006.6086

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H34").Value = 3003.5078
$ws.Range("I34").Value = 2301.6086
$ws.Range("K34").Value = 2301.6086
$ws.Range("M34").Value = -2099.6086

$ws = $wb.Sheets.Item("CUL")
$ws.Range("H2").Value = 976.6667
$ws.Range("I2").Value = 799.2308
$ws.Range("K2").Value = 4795.3848
$ws.Range("M2").Value = -4682.3848

$ws = $wb.Sheets.Item("CUL")
$ws.Range("H80").Value = 3672
$ws.Range("I80").Value = 2133.5
$ws.Range("K80").Value = 6400.5
$ws.Range("M80").Value = -5464.5

$ws = $wb.Sheets.Item("CUL")
$ws.Range("H83").Value = 3672
$ws.Range("I83").Value = 2133.5
$ws.Range("K83").Value = 19201.5
$ws.Range("M83").Value = -14521.5

$ws = $wb.Sheets.Item("CUL")
$ws.Range("H86").Value = 801.3333
$ws.Range("I86").Value = 801.3333
$ws.Range("K86").Value = 2403.9999
$ws.Range("M86").Value = -1217.9999

$ws = $wb.Sheets.Item("CUL")
$ws.Range("H89").Value = 801.3333
$ws.Range("I89").Value = 801.3333
$ws.Range("K89").Value = 7211.9997
$ws.Range("M89").Value = -1283.9997

$ws = $wb.Sheets.Item("GSM")
$ws.Range("H70").Value = 10677.131
$ws.Range("I70").Value = 10354.75
$ws.Range("K70").Value = 10354.75
$ws.Range("M70").Value = -10084.75

$ws = $wb.Sheets.Item("GSM")
$ws.Range("H73").Value = 10677.131
$ws.Range("I73").Value = 10354.75
$ws.Range("K73").Value = 10354.75
$ws.Range("M73").Value = -9418.75

$ws = $wb.Sheets.Item("GSM")
$ws.Range("H98").Value = 25045.092
$ws.Range("J98").Value = 25045.092
$ws.Range("L98").Value = 25045.092
$ws.Range("N98").Value = -31035.092

$ws = $wb.Sheets.Item("GSM")
$ws.Range("H122").Value = 2537.8
$ws.Range("I122").Value = 2312.8462
$ws.Range("K122").Value = 6938.5386
$ws.Range("M122").Value = -4488.5386

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H22").Value = 2481892.2
$ws.Range("I22").Value = 521.125
$ws.Range("J22").Value = 6452086
$ws.Range("K22").Value = 521.125
$ws.Range("L22").Value = 6452086
$ws.Range("M22").Value = -226.125
$ws.Range("N22").Value = -6452676

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H27").Value = 2481892.2
$ws.Range("I27").Value = 521.125
$ws.Range("J27").Value = 6452086
$ws.Range("K27").Value = 521.125
$ws.Range("L27").Value = 6452086
$ws.Range("M27").Value = -414.125
$ws.Range("N27").Value = -6452300

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H46").Value = 1906
$ws.Range("I46").Value = 1799.8334
$ws.Range("J46").Value = 2079.7273
$ws.Range("K46").Value = 1799.8334
$ws.Range("L46").Value = 2079.7273
$ws.Range("M46").Value = -1611.8334
$ws.Range("N46").Value = -2455.7273

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H82").Value = 2184.5925
$ws.Range("I82").Value = 2206.9443
$ws.Range("J82").Value = 2139.889
$ws.Range("K82").Value = 2206.9443
$ws.Range("L82").Value = 2139.889
$ws.Range("M82").Value = -1845.9443
$ws.Range("N82").Value = -2861.889

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H85").Value = 2184.5925
$ws.Range("I85").Value = 2206.9443
$ws.Range("J85").Value = 2139.889
$ws.Range("K85").Value = 2206.9443
$ws.Range("L85").Value = 2139.889
$ws.Range("M85").Value = -958.9443000000001
$ws.Range("N85").Value = -4635.889

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H93").Value = 4757.143
$ws.Range("I93").Value = 4500
$ws.Range("J93").Value = 4860
$ws.Range("K93").Value = 4500
$ws.Range("L93").Value = 4860
$ws.Range("M93").Value = -3252
$ws.Range("N93").Value = -7356

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H122").Value = 4792.143
$ws.Range("I122").Value = 3096.2354
$ws.Range("K122").Value = 9288.706200000001
$ws.Range("M122").Value = -6838.706200000001

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H132").Value = 5331.409
$ws.Range("J132").Value = 7886.222
$ws.Range("L132").Value = 23658.666
$ws.Range("N132").Value = -28718.666

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H126").Value = 1423.0769
$ws.Range("I126").Value = 1263.8695
$ws.Range("J126").Value = 1651.9375
$ws.Range("K126").Value = 3791.6085
$ws.Range("L126").Value = 4955.8125
$ws.Range("M126").Value = -1321.6085
$ws.Range("N126").Value = -9895.8125

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H132").Value = 5372.4287
$ws.Range("I132").Value = 4457.75
$ws.Range("J132").Value = 8299.4
$ws.Range("K132").Value = 13373.25
$ws.Range("L132").Value = 24898.2
$ws.Range("M132").Value = -10843.25
$ws.Range("N132").Value = -29958.2

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H136").Value = 2559.74
$ws.Range("I136").Value = 1886.6666
$ws.Range("K136").Value = 5659.9998
$ws.Range("M136").Value = -3109.9998
